$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column layout fix: the original "A:B" column-width definition overlapped
#    column B's own (wider) definition, so column B was rendered at column
#    A's width. Re-assert column B's width from column C (same intended
#    width) so the engine splits column A's range down to A:A only and
#    column B gets its correct width.
# ---------------------------------------------------------------------------
$ws.Range("B1").EntireColumn.ColumnWidth = $ws.Range("C1").EntireColumn.ColumnWidth

# ---------------------------------------------------------------------------
# 2) Row 10 (Objetivos:) - fill in the missing Portuguese objectives text
#    that had been left showing the professor's name by mistake.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value2 = "Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos e arranjos produtivos desta indústria."
$ws.Range("C10").Value2 = "Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos e arranjos produtivos desta indústria."

# ---------------------------------------------------------------------------
# 3) Row 12/13: insert the "Docentes responsáveis:" value ("... Barcza") as
#    its own row under the (already-present) label in A12, and remove the
#    stray "Programa resumido:" label that had incorrectly been placed
#    beside it.
# ---------------------------------------------------------------------------
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = "1285870 - Marcos Villela Barcza"
$ws.Range("C13").Value2 = "1285870 - Marcos Villela Barcza"

# ---------------------------------------------------------------------------
# 4) Row 14: "Programa resumido:" label + its real value (was blank/"Semestral").
# ---------------------------------------------------------------------------
$ws.Range("A14").Value2 = "Programa resumido:"
$ws.Range("B14").Value2 = "Processos da Indústria Químicos."
$ws.Range("C14").Value2 = "Processos da Indústria Químicos."

# ---------------------------------------------------------------------------
# 5) Row 15: "Short syllabus:" label + its real value (was showing a date).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value2 = "Short syllabus:"
$ws.Range("B15").Value2 = "Current and relevant topics related to chemical processes."
$ws.Range("C15").Value2 = "Current and relevant topics related to chemical processes."

# ---------------------------------------------------------------------------
# 6) Row 16: "Programa:" label + its real (Portuguese) syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "Programa:"
$ws.Range("B16").Value2 = "Panorama da Indústria Química. Química Fina. Petróleo e Petroquímica. Cerâmica. Vidro. Cimento. Celulose e Papel. Plásticos e afins. Processos Químicos Inovadores."
$ws.Range("C16").Value2 = "Panorama da Indústria Química. Química Fina. Petróleo e Petroquímica. Cerâmica. Vidro. Cimento. Celulose e Papel. Plásticos e afins. Processos Químicos Inovadores."

# ---------------------------------------------------------------------------
# 7) Row 17: "Syllabus:" label moves here, together with the English
#    syllabus text that used to sit on what is now row 16.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value2 = "Syllabus:"
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B17").Value2 = "Overview of the Chemical Industry. Fine Chemistry. Petroleum and Petrochemicals. Ceramics. Glass. Cement. Cellulose and paper. Plastics and allied products. Innovative Chemical Processes."
$ws.Range("C17").Value2 = "Overview of the Chemical Industry. Fine Chemistry. Petroleum and Petrochemicals. Ceramics. Glass. Cement. Cellulose and paper. Plastics and allied products. Innovative Chemical Processes."

# ---------------------------------------------------------------------------
# 8) Row 18: "Avaliação:" label moves here; the stray Barcza value that had
#    wrongly been duplicated under "Método:" is removed (real value now
#    lives on row 13).
# ---------------------------------------------------------------------------
$ws.Range("A18").Value2 = "Avaliação:"
$ws.Range("B18:C18").Clear()

# ---------------------------------------------------------------------------
# 9) Rows 19-21: labels shift down by one; their values were already correct
#    and stay in place.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value2 = "Método:"
$ws.Range("A20").Value2 = "Critério:"
$ws.Range("A21").Value2 = "Norma de recuperação:"

# ---------------------------------------------------------------------------
# 10) Row 22: "Bibliografia:" label moves here, together with its real
#     (previously missing) bibliography text.
# ---------------------------------------------------------------------------
$ws.Range("A22").Value2 = "Bibliografia:"
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B22").Value2 = "Material elaborado pelo docente.   LIVROS:Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim; New York: WileyVCH, 2011.  Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York: Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.   Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.   REVISTAS:Química & Derivados. Disponível em: http://www.quimica.com.br/pquimica/category/revista/Petróleo & Energia. Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
$ws.Range("C22").Value2 = "Material elaborado pelo docente.   LIVROS:Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim; New York: WileyVCH, 2011.  Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York: Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.   Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.   REVISTAS:Química & Derivados. Disponível em: http://www.quimica.com.br/pquimica/category/revista/Petróleo & Energia. Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."

# ---------------------------------------------------------------------------
# 11) Row 24 (new): the prerequisite text, moved one row down. Copy the
#     B23/C23 formatting down to B24/C24 BEFORE row 23's own content is
#     touched/cleared below.
# ---------------------------------------------------------------------------
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B24").Value2 = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Range("C24").Value2 = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"

# ---------------------------------------------------------------------------
# 12) Row 23: "Requisitos:" label moves here (was on row 22); the
#     prerequisite text that used to sit on this row is removed (it now
#     lives on the new row 24, filled in above).
# ---------------------------------------------------------------------------
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A23").Value2 = "Requisitos:"
$ws.Range("B23:C23").Clear()

# ---------------------------------------------------------------------------
# 13) Row heights - match the corrected row content/sizing.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(24).RowHeight = 30
